$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - sheet index 1
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value = 3137
$ws1.Range("F3").Value = 531
$ws1.Range("F4").Value = 1094
$ws1.Range("G5").Value = 35
$ws1.Range("F6").Value = 36
$ws1.Range("F8").Value = 38
$ws1.Range("F9").Value = 1127
$ws1.Range("F10").Value = 15708
$ws1.Range("F11").Value = 242
$ws1.Range("F12").Value = 178
$ws1.Range("F13").Value = 1025
$ws1.Range("F14").Value = 6182
$ws1.Range("F15").Value = 622
$ws1.Range("F16").Value = 108
$ws1.Range("F17").Value = 67
$ws1.Range("F18").Value = 7
$ws1.Range("F19").Value = 117
$ws1.Range("F20").Value = 1263
$ws1.Range("F22").Value = 633
$ws1.Range("F23").Value = 15
$ws1.Range("F26").Value = 210
$ws1.Range("F27").Value = 863
$ws1.Range("F28").Value = 29
$ws1.Range("F29").Value = 5001
$ws1.Range("F30").Value = 487
$ws1.Range("F31").Value = 11068
$ws1.Range("F32").Value = 1229
$ws1.Range("F33").Value = 13
$ws1.Range("F34").Value = 119
$ws1.Range("F35").Value = 169
$ws1.Range("F36").Value = 3804
$ws1.Range("F37").Value = 266

# Sheet "全部类型" (All types) - sheet index 4
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F3").Value = 3137
$ws4.Range("F4").Value = 531
$ws4.Range("F5").Value = 1094
$ws4.Range("G6").Value = 35
$ws4.Range("F7").Value = 36
$ws4.Range("F9").Value = 38
$ws4.Range("F10").Value = 1127
$ws4.Range("F11").Value = 15708
$ws4.Range("F12").Value = 242
$ws4.Range("F13").Value = 178
$ws4.Range("F14").Value = 1025
$ws4.Range("F15").Value = 6182
$ws4.Range("F16").Value = 622
$ws4.Range("F17").Value = 108
$ws4.Range("F18").Value = 67
$ws4.Range("F19").Value = 7
$ws4.Range("F20").Value = 117
$ws4.Range("F21").Value = 1263
$ws4.Range("F23").Value = 633
$ws4.Range("F24").Value = 15
$ws4.Range("F27").Value = 210
$ws4.Range("F28").Value = 863
$ws4.Range("F29").Value = 29
$ws4.Range("F30").Value = 5001
$ws4.Range("F31").Value = 487
$ws4.Range("F33").Value = 11068
$ws4.Range("F34").Value = 1229
$ws4.Range("F35").Value = 13
$ws4.Range("F36").Value = 119
$ws4.Range("F37").Value = 169
$ws4.Range("F38").Value = 3804
$ws4.Range("F39").Value = 266
